$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 867
$ws1.Range("F3").Value = 13845
$ws1.Range("F4").Value = 13625
$ws1.Range("F8").Value = 604
$ws1.Range("F12").Value = 770
$ws1.Range("F13").Value = 2152
$ws1.Range("F14").Value = 113
$ws1.Range("G14").Value = 68
$ws1.Range("F16").Value = 78
$ws1.Range("G16").Value = 60
$ws1.Range("F19").Value = 536
$ws1.Range("F20").Value = 437
$ws1.Range("F21").Value = 414
$ws1.Range("F23").Value = 271

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 49
$ws2.Range("F7").Value = 1539

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 114

# --- Sheet "全部类型" (All types, merged view) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 867
$ws4.Range("F4").Value = 13845
$ws4.Range("F5").Value = 13625
$ws4.Range("F9").Value = 604
$ws4.Range("F13").Value = 770
$ws4.Range("F14").Value = 49
$ws4.Range("F16").Value = 2152
$ws4.Range("F17").Value = 113
$ws4.Range("G17").Value = 68
$ws4.Range("F19").Value = 78
$ws4.Range("G19").Value = 60
$ws4.Range("F24").Value = 114
$ws4.Range("F25").Value = 114
$ws4.Range("F26").Value = 536
$ws4.Range("F27").Value = 437
$ws4.Range("F28").Value = 414
$ws4.Range("F30").Value = 271
$ws4.Range("F33").Value = 1539
